$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in column D for rows 2-5:
# Rows 2-3 had 44574 -> now 44559
# Rows 4-5 had 44559 -> now 44574
$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D4").Value = 44574
$ws.Range("D5").Value = 44574
